$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# Insert two new blank columns before the existing data column C.
# This pushes the old "B" (week header / UN placeholder) column to D
# and the old "C" (actual rank data) column to E, exactly like the
# target workbook's layout.
# ------------------------------------------------------------------
$ws.Columns("B:C").Insert()

# ------------------------------------------------------------------
# Row 1 (headers): new columns get the two new week labels.
# D1/E1 already hold the old header values (Jun_13 / Jun_10) because
# of the column insert above.
# ------------------------------------------------------------------
$ws.Range("B1").Value = "Jun_17"
$ws.Range("C1").Value = "Jun_15"

# ------------------------------------------------------------------
# Data rows: fill the two new columns with the same "UN" placeholder
# used throughout column B previously (old default value), mirroring
# the pattern applied to the rest of the sheet.
# ------------------------------------------------------------------
for ($r = 2; $r -le 27; $r++) {
    $ws.Cells.Item($r, 2).Value = "UN"
    $ws.Cells.Item($r, 3).Value = "UN"
}

# ------------------------------------------------------------------
# Cosmetic column widths: make the three "week data" columns (C, D, E)
# share the same custom width that the original data column (C) used
# to have on its own.
# ------------------------------------------------------------------
$ws.Columns("C").ColumnWidth = 7.1
$ws.Columns("D").ColumnWidth = 7.1
$ws.Columns("E").ColumnWidth = 7.1
